$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting the existing rows 250-379 down
# to 251-380 (dimension grows from A1:R379 to A1:R380).
$ws.Rows("250:250").Insert()

# Populate the newly-inserted row 250 with its data.
$ws.Range("A250").Value = 9
$ws.Range("B250").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C250").Value = 'Metropolitana'
$ws.Range("D250").Value = 44960
$ws.Range("E250").Value = 13
$ws.Range("F250").Value = 300000001
$ws.Range("G250").Value = 'Rabanito'
$ws.Range("H250").Value = 'Sin especificar'
$ws.Range("I250").Value = 'Primera'
$ws.Range("J250").Value = 7000
$ws.Range("K250").Value = 3000
$ws.Range("L250").Value = 3000
$ws.Range("M250").Value = 3000
$ws.Range("N250").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O250").Value = 'Provincia de Chacabuco'
$ws.Range("P250").Value = 30
$ws.Range("Q250").Value = 100
$ws.Range("R250").Value = 'Hortaliza'
